# "majors to job csvs"
# Applies the OOXML diff to The FortunED Tellers.pptx via PowerPoint COM automation.

$p = $ppt.ActivePresentation

# ---------------------------------------------------------------------------
# 1) Footer date placeholders on every slide layout + the slide master:
#    "6/30/2020" -> "7/1/2020"
# ---------------------------------------------------------------------------
$master = $p.SlideMaster
$masterDateShape = $null
foreach ($shp in $master.Shapes) {
    if ($shp.HasTextFrame -and $shp.TextFrame.HasText) {
        if ($shp.TextFrame.TextRange.Text -eq "6/30/2020") {
            $shp.TextFrame.TextRange.Text = "7/1/2020"
        }
    }
}

foreach ($layout in $master.CustomLayouts) {
    foreach ($shp in $layout.Shapes) {
        if ($shp.HasTextFrame -and $shp.TextFrame.HasText) {
            if ($shp.TextFrame.TextRange.Text -eq "6/30/2020") {
                $shp.TextFrame.TextRange.Text = "7/1/2020"
            }
        }
    }
}

# ---------------------------------------------------------------------------
# 2) Slide 3 ("Page 2 - HS") content placeholder
# ---------------------------------------------------------------------------
$s3 = $p.Slides.Item(3)
$tr3 = $s3.Shapes.Item(2).TextFrame.TextRange

# Remove trailing empty paragraph (lvl 1) at the end of the body.
$tr3.Paragraphs(11).Delete()

$tr3.Paragraphs(10).Runs(1).Text = "Career options within Industry"
$tr3.Paragraphs(9).Runs(1).Text = "Starting – Mid – Senior Salary Projections"

# ---------------------------------------------------------------------------
# 3) Slide 4 ("Page 2 - College") content placeholder
# ---------------------------------------------------------------------------
$s4 = $p.Slides.Item(4)
$shp4 = $s4.Shapes.Item(2)
$tr4 = $shp4.TextFrame.TextRange

# Work from the bottom of the list upwards so paragraph indices of
# not-yet-processed items stay stable.

# Para 12: "Scrape indeed for Internships..." -> add orange font color (FFC000)
$tr4.Paragraphs(12).Font.Color.RGB = 49407

# Para 11: "How much they have to pay every month to pay off student loan (ML)"
#   -> "How long it will take to pay off student loan (ML)"
#   + new lvl2 paragraph after it.
$tr4.Paragraphs(11).Runs(1).Text = "How long it will take to pay off student loan (ML)"
$newP = $tr4.Paragraphs(11).InsertAfter([char]13 + "Tuition, Starting Salary, Cost of living, Question: What percent of salary will you spend towards Loan?")
$tr4.Paragraphs(12).IndentLevel = 3

# Para 10: "Rank of cities to job hunt" (lvl2)
#   -> "Options of cities more favorable based on starting salary" + orange color
$tr4.Paragraphs(10).Runs(1).Text = "Options of cities more favorable based on starting salary"
$tr4.Paragraphs(10).Font.Color.RGB = 49407

# Para 9: "Cost of living vs Starting Salary (create algorithm)" (lvl1)
#   -> "Cost of living vs P25 starting salary (create algorithm)" + orange color
#   + new lvl1 paragraph "Show range of starting salaries" inserted before it.
$tr4.Paragraphs(9).Runs(1).Text = "Cost of living vs P25 starting salary (create algorithm)"
$tr4.Paragraphs(9).Font.Color.RGB = 49407
$null = $tr4.Paragraphs(8).InsertAfter([char]13 + "Show range of starting salaries")

# Para 5 -> "How many years do you want to pay it off?"
#   + new lvl0 paragraph "Ideal City/State to work in?" inserted after it.
$null = $tr4.Paragraphs(5).InsertAfter([char]13 + "Ideal City/State to work in?")

# ---------------------------------------------------------------------------
# 4) Slide 5 ("Parent 2 - Parent/Guardian") content placeholder
# ---------------------------------------------------------------------------
$s5 = $p.Slides.Item(5)
$tr5 = $s5.Shapes.Item(2).TextFrame.TextRange
$full5 = $tr5.Text
$needle = "graphs from here: "
$idx = $full5.IndexOf($needle)
$sub = $tr5.Characters($idx + 1, $needle.Length)
$sub.Delete()
$tr5.Paragraphs(4).Runs(1).Text = "Some graphs from here: "
